{"js": "// Replace the date line and the 25 two-digit multiplication equations\n// with the updated values from the commit.\nconst replacements = [\n  [\"2023-09-26 Tuesday\", \"2023-09-27 Wednesday\"],\n  [\"82\u00d789=7298\", \"90\u00d743=3870\"],\n  [\"24\u00d773=1752\", \"93\u00d794=8742\"],\n  [\"22\u00d795=2090\", \"71\u00d745=3195\"],\n  [\"36\u00d722=792\", \"43\u00d789=3827\"],\n  [\"53\u00d715=795\", \"52\u00d713=676\"],\n  [\"67\u00d718=1206\", \"16\u00d751=816\"],\n  [\"46\u00d741=1886\", \"43\u00d717=731\"],\n  [\"17\u00d713=221\", \"53\u00d730=1590\"],\n  [\"60\u00d770=4200\", \"28\u00d752=1456\"],\n  [\"66\u00d776=5016\", \"24\u00d753=1272\"],\n  [\"74\u00d734=2516\", \"13\u00d737=481\"],\n  [\"35\u00d769=2415\", \"95\u00d782=7790\"],\n  [\"55\u00d752=2860\", \"21\u00d798=2058\"],\n  [\"63\u00d789=5607\", \"75\u00d793=6975\"],\n  [\"42\u00d755=2310\", \"60\u00d712=720\"],\n  [\"98\u00d732=3136\", \"20\u00d762=1240\"],\n  [\"80\u00d745=3600\", \"99\u00d784=8316\"],\n  [\"23\u00d760=1380\", \"18\u00d718=324\"],\n  [\"24\u00d734=816\", \"74\u00d791=6734\"],\n  [\"75\u00d742=3150\", \"17\u00d776=1292\"],\n  [\"95\u00d727=2565\", \"62\u00d758=3596\"],\n  [\"46\u00d758=2668\", \"55\u00d724=1320\"],\n  [\"16\u00d746=736\", \"11\u00d743=473\"],\n  [\"21\u00d756=1176\", \"47\u00d759=2773\"],\n  [\"82\u00d769=5658\", \"77\u00d713=1001\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 two-digit multiplication equations\n# with the updated values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-09-26 Tuesday\", \"2023-09-27 Wednesday\"),\n    @(\"82\u00d789=7298\", \"90\u00d743=3870\"),\n    @(\"24\u00d773=1752\", \"93\u00d794=8742\"),\n    @(\"22\u00d795=2090\", \"71\u00d745=3195\"),\n    @(\"36\u00d722=792\", \"43\u00d789=3827\"),\n    @(\"53\u00d715=795\", \"52\u00d713=676\"),\n    @(\"67\u00d718=1206\", \"16\u00d751=816\"),\n    @(\"46\u00d741=1886\", \"43\u00d717=731\"),\n    @(\"17\u00d713=221\", \"53\u00d730=1590\"),\n    @(\"60\u00d770=4200\", \"28\u00d752=1456\"),\n    @(\"66\u00d776=5016\", \"24\u00d753=1272\"),\n    @(\"74\u00d734=2516\", \"13\u00d737=481\"),\n    @(\"35\u00d769=2415\", \"95\u00d782=7790\"),\n    @(\"55\u00d752=2860\", \"21\u00d798=2058\"),\n    @(\"63\u00d789=5607\", \"75\u00d793=6975\"),\n    @(\"42\u00d755=2310\", \"60\u00d712=720\"),\n    @(\"98\u00d732=3136\", \"20\u00d762=1240\"),\n    @(\"80\u00d745=3600\", \"99\u00d784=8316\"),\n    @(\"23\u00d760=1380\", \"18\u00d718=324\"),\n    @(\"24\u00d734=816\", \"74\u00d791=6734\"),\n    @(\"75\u00d742=3150\", \"17\u00d776=1292\"),\n    @(\"95\u00d727=2565\", \"62\u00d758=3596\"),\n    @(\"46\u00d758=2668\", \"55\u00d724=1320\"),\n    @(\"16\u00d746=736\", \"11\u00d743=473\"),\n    @(\"21\u00d756=1176\", \"47\u00d759=2773\"),\n    @(\"82\u00d769=5658\", \"77\u00d713=1001\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
